$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9994668960571289
$ws.Range("B1").Value = 1.024011135101318
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.095846176147461
$ws.Range("E1").Value = 1.053468942642212
